$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add new header + value for the Zipcode column (D)
$ws.Range("D1").Value = "Zipcode"
$ws.Range("D2").Value = 500062

# Select D2 like the source document shows
$ws.Range("D2").Select()
